$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2199.2
$ws.Range("I40").Value = 1998.6666
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 1998.6666
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -1823.6666
$ws.Range("N40").Value = -2850
# Row 62
$ws.Range("H62").Value = 4160.125
$ws.Range("I62").Value = 3770.8
$ws.Range("J62").Value = 10000
$ws.Range("K62").Value = 3770.8
$ws.Range("L62").Value = 10000
$ws.Range("M62").Value = -3146.8
$ws.Range("N62").Value = -11248
# Row 65
$ws.Range("H65").Value = 4160.125
$ws.Range("I65").Value = 3770.8
$ws.Range("J65").Value = 10000
$ws.Range("K65").Value = 18854
$ws.Range("L65").Value = 50000
$ws.Range("M65").Value = -15734
$ws.Range("N65").Value = -56240
# Row 135
$ws.Range("H135").Value = 1537.0588
$ws.Range("I135").Value = 1289.6364
$ws.Range("J135").Value = 1990.6666
$ws.Range("K135").Value = 11606.7276
$ws.Range("L135").Value = 17915.9994
$ws.Range("M135").Value = -9071.7276
$ws.Range("N135").Value = -22985.9994
# Row 137
$ws.Range("H137").Value = 1640.5
$ws.Range("I137").Value = 1581.7778
$ws.Range("J137").Value = 1816.6666
$ws.Range("K137").Value = 4745.3334
$ws.Range("L137").Value = 5449.9998
$ws.Range("M137").Value = -2195.3334
$ws.Range("N137").Value = -10549.9998

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4191.091
$ws.Range("I32").Value = 2668.4878
$ws.Range("J32").Value = 25000
$ws.Range("K32").Value = 2668.4878
$ws.Range("L32").Value = 25000
$ws.Range("M32").Value = -2381.4878
$ws.Range("N32").Value = -25574
# Row 45
$ws.Range("H45").Value = 1999.5
$ws.Range("I45").Value = 1999.5
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1999.5
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -1622.5
# Row 48
$ws.Range("H48").Value = 170000
$ws.Range("J48").Value = 170000
$ws.Range("L48").Value = 170000
$ws.Range("N48").Value = -170768
# Row 63
$ws.Range("H63").Value = 3196.4285
$ws.Range("I63").Value = 2093.75
$ws.Range("J63").Value = 4666.6665
$ws.Range("K63").Value = 2093.75
$ws.Range("L63").Value = 4666.6665
$ws.Range("M63").Value = -1407.75
$ws.Range("N63").Value = -6038.6665
# Row 66
$ws.Range("H66").Value = 3196.4285
$ws.Range("I66").Value = 2093.75
$ws.Range("J66").Value = 4666.6665
$ws.Range("K66").Value = 10468.75
$ws.Range("L66").Value = 23333.3325
$ws.Range("M66").Value = -7036.75
$ws.Range("N66").Value = -30197.3325
# Row 123
$ws.Range("H123").Value = 66666.664
$ws.Range("J123").Value = 66666.664
$ws.Range("L123").Value = 66666.664
$ws.Range("N123").Value = -76466.664

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 41
$ws.Range("H41").Value = 170000
$ws.Range("J41").Value = 170000
$ws.Range("L41").Value = 170000
$ws.Range("N41").Value = -170776
# Row 47
$ws.Range("H47").Value = 170000
$ws.Range("J47").Value = 170000
$ws.Range("L47").Value = 170000
$ws.Range("N47").Value = -171040
# Row 64
$ws.Range("H64").Value = 1502
$ws.Range("I64").Value = 500
$ws.Range("K64").Value = 500
$ws.Range("M64").Value = -275
# Row 67
$ws.Range("H67").Value = 1502
$ws.Range("I67").Value = 500
$ws.Range("K67").Value = 500
$ws.Range("M67").Value = 280

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4614.4707
$ws.Range("I31").Value = 3433.6365
$ws.Range("J31").Value = 6779.3335
$ws.Range("K31").Value = 3433.6365
$ws.Range("L31").Value = 6779.3335
$ws.Range("M31").Value = -3138.6365
$ws.Range("N31").Value = -7369.3335
# Row 32
$ws.Range("H32").Value = 354.54544
$ws.Range("I32").Value = 354.54544
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 354.54544
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -38.54543999999999
# Row 34
$ws.Range("H34").Value = 4614.4707
$ws.Range("I34").Value = 3433.6365
$ws.Range("J34").Value = 6779.3335
$ws.Range("K34").Value = 3433.6365
$ws.Range("L34").Value = 6779.3335
$ws.Range("M34").Value = -3231.6365
$ws.Range("N34").Value = -7183.3335
# Row 97
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").ClearContents()
$ws.Range("N97").Value = 0

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 2133
$ws.Range("J5").Value = 2748
$ws.Range("L5").Value = 8244
$ws.Range("N5").Value = -8468
# Row 47
$ws.Range("H47").Value = 53
$ws.Range("I47").Value = 53
$ws.Range("K47").Value = 159
$ws.Range("M47").Value = 272
# Row 60
$ws.Range("H60").Value = 503.5
$ws.Range("I60").Value = 138
$ws.Range("J60").Value = 1600
$ws.Range("K60").Value = 414
$ws.Range("L60").Value = 4800
$ws.Range("M60").Value = -163
$ws.Range("N60").Value = -5302
# Row 118
$ws.Range("H118").Value = 417.16666
$ws.Range("I118").Value = 417.16666
$ws.Range("K118").Value = 1251.49998
$ws.Range("M118").Value = -8.49998000000005
# Row 120
$ws.Range("H120").Value = 9979.866
$ws.Range("I120").Value = 5966.4443
$ws.Range("K120").Value = 17899.3329
$ws.Range("M120").Value = -13061.3329
# Row 131
$ws.Range("H131").Value = 1062.4
$ws.Range("J131").Value = 1384.2307
$ws.Range("L131").Value = 4152.6921
$ws.Range("N131").Value = -14232.6921
# Row 135
$ws.Range("H135").Value = 2133
$ws.Range("J135").Value = 2748
$ws.Range("L135").Value = 24732
$ws.Range("N135").Value = -29802

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 26
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").ClearContents()
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = 0
# Row 50
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").ClearContents()
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = 0
# Row 97
$ws.Range("H97").Value = 1216.85
$ws.Range("I97").Value = 1214.7
$ws.Range("J97").Value = 1219
$ws.Range("K97").Value = 1214.7
$ws.Range("L97").Value = 1219
$ws.Range("M97").Value = -718.7
$ws.Range("N97").Value = -2211
# Row 126
$ws.Range("H126").Value = 3701.375
$ws.Range("I126").Value = 3474.5
$ws.Range("J126").Value = 4382
$ws.Range("K126").Value = 10423.5
$ws.Range("L126").Value = 13146
$ws.Range("M126").Value = -7953.5
$ws.Range("N126").Value = -18086
# Row 136
$ws.Range("H136").Value = 23285.143
$ws.Range("J136").Value = 23285.143
$ws.Range("L136").Value = 69855.429
$ws.Range("N136").Value = -74955.429

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 6305.697
$ws.Range("I22").Value = 5082.933
$ws.Range("J22").Value = 7324.6665
$ws.Range("K22").Value = 5082.933
$ws.Range("L22").Value = 7324.6665
$ws.Range("M22").Value = -4787.933
$ws.Range("N22").Value = -7914.6665
# Row 27
$ws.Range("H27").Value = 6305.697
$ws.Range("I27").Value = 5082.933
$ws.Range("J27").Value = 7324.6665
$ws.Range("K27").Value = 5082.933
$ws.Range("L27").Value = 7324.6665
$ws.Range("M27").Value = -4975.933
$ws.Range("N27").Value = -7538.6665
# Row 46
$ws.Range("H46").Value = 3018.4814
$ws.Range("I46").Value = 2000.9231
$ws.Range("J46").Value = 3963.3572
$ws.Range("K46").Value = 2000.9231
$ws.Range("L46").Value = 3963.3572
$ws.Range("M46").Value = -1812.9231
$ws.Range("N46").Value = -4339.3572
# Row 50
$ws.Range("H50").Value = 60084
$ws.Range("J50").Value = 60084
$ws.Range("L50").Value = 60084
$ws.Range("N50").Value = -61358
# Row 93
$ws.Range("H93").Value = 1550.375
$ws.Range("I93").Value = 2100
$ws.Range("J93").Value = 1000.75
$ws.Range("K93").Value = 2100
$ws.Range("L93").Value = 1000.75
$ws.Range("M93").Value = -852
$ws.Range("N93").Value = -3496.75
# Row 125
$ws.Range("H125").Value = 75000
$ws.Range("J125").Value = 75000
$ws.Range("L125").Value = 75000
$ws.Range("N125").Value = -84840
# Row 132
$ws.Range("H132").Value = 5232.067
$ws.Range("J132").Value = 6070.5713
$ws.Range("L132").Value = 18211.7139
$ws.Range("N132").Value = -23271.7139

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 21
$ws.Range("H21").Value = 19375
# Row 24
$ws.Range("H24").Value = 30010
$ws.Range("J24").Value = 30010
$ws.Range("L24").Value = 30010
$ws.Range("N24").Value = -30470
# Row 35
$ws.Range("H35").Value = 19375
# Row 49
$ws.Range("H49").Value = 500028
$ws.Range("I49").Value = 500028
$ws.Range("K49").Value = 500028
$ws.Range("M49").Value = -499798
